# Apply the "Updated symbol list" edit: update Price (column D) values and a
# handful of combined-label (column E) strings for specific rows.
#
# The Price column cells are stored as text (e.g. "244.50"), so we force the
# cell's number format to Text ("@") before assigning the new value. This
# keeps values such as "3.250" or "0.00000000750" exactly as written instead
# of Excel silently re-interpreting them as numbers and dropping trailing
# zeros / significant digits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Column D (Price) updates
Set-TextCell "D2"  "244.50"
Set-TextCell "D4"  "5.747"
Set-TextCell "D5"  "0.05831"
Set-TextCell "D6"  "3.415"
Set-TextCell "D7"  "6.469"
Set-TextCell "D8"  "1.320"
Set-TextCell "D9"  "0.8023"
Set-TextCell "D10" "0.1464"
Set-TextCell "D12" "0.03251"
Set-TextCell "D13" "0.03003"
Set-TextCell "D14" "0.09247"
Set-TextCell "D15" "0.001661"
Set-TextCell "D16" "3.252"
Set-TextCell "D17" "0.04769"
Set-TextCell "D18" "0.0005988"
Set-TextCell "D19" "0.006255"
Set-TextCell "D20" "0.005422"
Set-TextCell "D21" "0.001065"
Set-TextCell "D22" "0.0001499"
Set-TextCell "D23" "3.697"
Set-TextCell "D26" "0.1267"
Set-TextCell "D27" "0.0009994"
Set-TextCell "D41" "0.007054"
Set-TextCell "D42" "0.1062"
Set-TextCell "D43" "0.003168"
Set-TextCell "D44" "0.009763"
Set-TextCell "D46" "0.00005602"
Set-TextCell "D47" "0.00000000750"
Set-TextCell "D48" "0.7849"
Set-TextCell "D49" "0.09894"
Set-TextCell "D50" "0.00002099"
Set-TextCell "D51" "0.01010"

# Column E (Coin+Symbol+"Best/Worst in24h" label) updates
Set-TextCell "E18" "17OneONEWorstin24h"
Set-TextCell "E27" "26UpBotsUBXTBestin24h"
Set-TextCell "E41" "40KickTokenKICK"
Set-TextCell "E45" "44ACDXExchangeACXT"
